$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 170.16317131227225
$ws.Range("C2").Value = 140.01849556838636
$ws.Range("D2").Value = 171.0798785168999
$ws.Range("E2").Value = 136.13812285768728

$ws.Range("B3").Value = 153.98849238367626
$ws.Range("C3").Value = 130.37329613799668
$ws.Range("D3").Value = 160.1162219257493
$ws.Range("E3").Value = 132.86164789905024

$ws.Range("B1:E3").Select()
